$d = $word.ActiveDocument

# The document ends with:
#   ... bibliography paragraph ending in "3040P."
#   <empty paragraph>
#   "Ver no Jupiter Salvar em pdf Salvar em docx"
#   "(c) 2020 . Contact: ... Original theme under Creative Commons Attribution"
#   <empty paragraph>
#   <page-break paragraph>
#
# We need to remove the empty paragraph, the "Ver no Jupiter..." paragraph
# and the "(c) 2020..." paragraph, leaving the bibliography paragraph
# directly followed by the trailing empty paragraph and the page break.

$footerStartText = "Ver no Jupiter Salvar em pdf Salvar em docx"
$footerEndText = "Powered by Jekyll and Github pages"

$startIndex = -1
$endIndex = -1

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $paraText = $d.Paragraphs.Item($i).Range.Text
    if ($paraText -like "*$footerStartText*") {
        $startIndex = $i
    }
    if ($paraText -like "*$footerEndText*") {
        $endIndex = $i
    }
}

if ($startIndex -gt 0 -and $endIndex -ge $startIndex) {
    # Also remove the blank paragraph immediately preceding the
    # "Ver no Jupiter..." paragraph.
    $removeFrom = $startIndex
    $precedingIndex = $startIndex - 1
    if ($precedingIndex -ge 1) {
        $precedingText = $d.Paragraphs.Item($precedingIndex).Range.Text
        if ($precedingText.Trim().Length -eq 0) {
            $removeFrom = $precedingIndex
        }
    }

    $rangeStart = $d.Paragraphs.Item($removeFrom).Range.Start
    $rangeEnd = $d.Paragraphs.Item($endIndex).Range.End
    $d.Range($rangeStart, $rangeEnd).Delete()
}
